$wb = $excel.ActiveWorkbook

# ---- Step 1: remove the original "总计" (totals) sheet ----
$totalOld = $wb.Worksheets.Item(6)
$totalOld.Delete()

# ---- Step 2: build "2022-Q1" by copying "2021-Q4" (same column layout/style) ----
$q4 = $wb.Worksheets.Item(5)
$q4.Copy($null, $q4)
$newQ1 = $wb.Worksheets.Item(6)
$newQ1.Name = "2022-Q1"

# ---- Step 3: build the new "总计" by copying "2021-Q4" again (keeps the same bold/bordered style) ----
$q4b = $wb.Worksheets.Item(5)
$q4b.Copy($null, $newQ1)
$newTotal = $wb.Worksheets.Item(7)
$newTotal.Name = "总计"

# ================= Fill "2022-Q1" sheet =================
$ws = $newQ1
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

$ws.Range("A2").Value = 0
$c = $ws.Range("B2")
$c.NumberFormat = "@"
$c.Value = "005669"
$c.Style = "Normal"
$c = $ws.Range("C2")
$c.NumberFormat = "@"
$c.Value = "前海开源公用事业行业股票"
$c.Style = "Normal"
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "258.16"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "94.53"
$c.Style = "Normal"
$c = $ws.Range("F2")
$c.NumberFormat = "@"
$c.Value = "5.69"
$c.Style = "Normal"
$c = $ws.Range("G2")
$c.NumberFormat = "@"
$c.Value = "14.6893"
$c.Style = "Normal"
$ws.Range("H2").Value = 6
$ws.Range("A3").Value = 1
$c = $ws.Range("B3")
$c.NumberFormat = "@"
$c.Value = "001875"
$c.Style = "Normal"
$c = $ws.Range("C3")
$c.NumberFormat = "@"
$c.Value = "前海开源沪港深优势精选灵活配置混合"
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "82.95"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "93.80"
$c.Style = "Normal"
$c = $ws.Range("F3")
$c.NumberFormat = "@"
$c.Value = "4.91"
$c.Style = "Normal"
$c = $ws.Range("G3")
$c.NumberFormat = "@"
$c.Value = "4.0728"
$c.Style = "Normal"
$ws.Range("H3").Value = 9
$ws.Range("A4").Value = 2
$c = $ws.Range("B4")
$c.NumberFormat = "@"
$c.Value = "001874"
$c.Style = "Normal"
$c = $ws.Range("C4")
$c.NumberFormat = "@"
$c.Value = "前海开源沪港深价值精选灵活配置混合"
$c.Style = "Normal"
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "9.62"
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "94.37"
$c.Style = "Normal"
$c = $ws.Range("F4")
$c.NumberFormat = "@"
$c.Value = "5.06"
$c.Style = "Normal"
$c = $ws.Range("G4")
$c.NumberFormat = "@"
$c.Value = "0.4868"
$c.Style = "Normal"
$ws.Range("H4").Value = 8
$ws.Range("A5").Value = 3
$c = $ws.Range("B5")
$c.NumberFormat = "@"
$c.Value = "013270"
$c.Style = "Normal"
$c = $ws.Range("C5")
$c.NumberFormat = "@"
$c.Value = "前海开源聚利一年持有混合A"
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "7.55"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "80.39"
$c.Style = "Normal"
$c = $ws.Range("F5")
$c.NumberFormat = "@"
$c.Value = "6.29"
$c.Style = "Normal"
$c = $ws.Range("G5")
$c.NumberFormat = "@"
$c.Value = "0.4749"
$c.Style = "Normal"
$ws.Range("H5").Value = 7
$ws.Range("A6").Value = 4
$c = $ws.Range("B6")
$c.NumberFormat = "@"
$c.Value = "010751"
$c.Style = "Normal"
$c = $ws.Range("C6")
$c.NumberFormat = "@"
$c.Value = "宝盈优质成长混合A"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "5.64"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "92.80"
$c.Style = "Normal"
$c = $ws.Range("F6")
$c.NumberFormat = "@"
$c.Value = "5.06"
$c.Style = "Normal"
$c = $ws.Range("G6")
$c.NumberFormat = "@"
$c.Value = "0.2854"
$c.Style = "Normal"
$ws.Range("H6").Value = 5
$ws.Range("A7").Value = 5
$c = $ws.Range("B7")
$c.NumberFormat = "@"
$c.Value = "002653"
$c.Style = "Normal"
$c = $ws.Range("C7")
$c.NumberFormat = "@"
$c.Value = "泰康沪港深精选灵活配置混合"
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "7.89"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "87.05"
$c.Style = "Normal"
$c = $ws.Range("F7")
$c.NumberFormat = "@"
$c.Value = "2.07"
$c.Style = "Normal"
$c = $ws.Range("G7")
$c.NumberFormat = "@"
$c.Value = "0.1633"
$c.Style = "Normal"
$ws.Range("H7").Value = 10
$ws.Range("A8").Value = 6
$c = $ws.Range("B8")
$c.NumberFormat = "@"
$c.Value = "011157"
$c.Style = "Normal"
$c = $ws.Range("C8")
$c.NumberFormat = "@"
$c.Value = "弘毅远方港股通智选领航混合A"
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "3.65"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "90.35"
$c.Style = "Normal"
$c = $ws.Range("F8")
$c.NumberFormat = "@"
$c.Value = "4.16"
$c.Style = "Normal"
$c = $ws.Range("G8")
$c.NumberFormat = "@"
$c.Value = "0.1518"
$c.Style = "Normal"
$ws.Range("H8").Value = 9
$ws.Range("A9").Value = 7
$c = $ws.Range("B9")
$c.NumberFormat = "@"
$c.Value = "008404"
$c.Style = "Normal"
$c = $ws.Range("C9")
$c.NumberFormat = "@"
$c.Value = "华泰紫金泰盈混合A"
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "4.29"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "79.73"
$c.Style = "Normal"
$c = $ws.Range("F9")
$c.NumberFormat = "@"
$c.Value = "3.32"
$c.Style = "Normal"
$c = $ws.Range("G9")
$c.NumberFormat = "@"
$c.Value = "0.1424"
$c.Style = "Normal"
$ws.Range("H9").Value = 10
$ws.Range("A10").Value = 8
$c = $ws.Range("B10")
$c.NumberFormat = "@"
$c.Value = "008405"
$c.Style = "Normal"
$c = $ws.Range("C10")
$c.NumberFormat = "@"
$c.Value = "华泰紫金泰盈混合C"
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "3.71"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "79.73"
$c.Style = "Normal"
$c = $ws.Range("F10")
$c.NumberFormat = "@"
$c.Value = "3.32"
$c.Style = "Normal"
$c = $ws.Range("G10")
$c.NumberFormat = "@"
$c.Value = "0.1232"
$c.Style = "Normal"
$ws.Range("H10").Value = 10
$ws.Range("A11").Value = 9
$c = $ws.Range("B11")
$c.NumberFormat = "@"
$c.Value = "011694"
$c.Style = "Normal"
$c = $ws.Range("C11")
$c.NumberFormat = "@"
$c.Value = "华泰紫金信息科技主题6个月定期开放混合A"
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "2.60"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "77.49"
$c.Style = "Normal"
$c = $ws.Range("F11")
$c.NumberFormat = "@"
$c.Value = "3.24"
$c.Style = "Normal"
$c = $ws.Range("G11")
$c.NumberFormat = "@"
$c.Value = "0.0842"
$c.Style = "Normal"
$ws.Range("H11").Value = 10
$ws.Range("A12").Value = 10
$c = $ws.Range("B12")
$c.NumberFormat = "@"
$c.Value = "011651"
$c.Style = "Normal"
$c = $ws.Range("C12")
$c.NumberFormat = "@"
$c.Value = "招商港股通核心精选股票A"
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "2.81"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "81.27"
$c.Style = "Normal"
$c = $ws.Range("F12")
$c.NumberFormat = "@"
$c.Value = "2.94"
$c.Style = "Normal"
$c = $ws.Range("G12")
$c.NumberFormat = "@"
$c.Value = "0.0826"
$c.Style = "Normal"
$ws.Range("H12").Value = 5
$ws.Range("A13").Value = 11
$c = $ws.Range("B13")
$c.NumberFormat = "@"
$c.Value = "003580"
$c.Style = "Normal"
$c = $ws.Range("C13")
$c.NumberFormat = "@"
$c.Value = "泰康沪港深价值优选灵活配置混合"
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.80"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "85.70"
$c.Style = "Normal"
$c = $ws.Range("F13")
$c.NumberFormat = "@"
$c.Value = "2.27"
$c.Style = "Normal"
$c = $ws.Range("G13")
$c.NumberFormat = "@"
$c.Value = "0.0409"
$c.Style = "Normal"
$ws.Range("H13").Value = 9
$ws.Range("A14").Value = 12
$c = $ws.Range("B14")
$c.NumberFormat = "@"
$c.Value = "010752"
$c.Style = "Normal"
$c = $ws.Range("C14")
$c.NumberFormat = "@"
$c.Value = "宝盈优质成长混合C"
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.78"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "92.80"
$c.Style = "Normal"
$c = $ws.Range("F14")
$c.NumberFormat = "@"
$c.Value = "5.06"
$c.Style = "Normal"
$c = $ws.Range("G14")
$c.NumberFormat = "@"
$c.Value = "0.0395"
$c.Style = "Normal"
$ws.Range("H14").Value = 5
$ws.Range("A15").Value = 13
$c = $ws.Range("B15")
$c.NumberFormat = "@"
$c.Value = "004266"
$c.Style = "Normal"
$c = $ws.Range("C15")
$c.NumberFormat = "@"
$c.Value = "招商沪港深科技创新主题精选灵活配置混合A"
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "1.29"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "88.85"
$c.Style = "Normal"
$c = $ws.Range("F15")
$c.NumberFormat = "@"
$c.Value = "2.94"
$c.Style = "Normal"
$c = $ws.Range("G15")
$c.NumberFormat = "@"
$c.Value = "0.0379"
$c.Style = "Normal"
$ws.Range("H15").Value = 8
$ws.Range("A16").Value = 14
$c = $ws.Range("B16")
$c.NumberFormat = "@"
$c.Value = "006923"
$c.Style = "Normal"
$c = $ws.Range("C16")
$c.NumberFormat = "@"
$c.Value = "前海开源沪港深非周期性行业股票A"
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.54"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "93.77"
$c.Style = "Normal"
$c = $ws.Range("F16")
$c.NumberFormat = "@"
$c.Value = "6.16"
$c.Style = "Normal"
$c = $ws.Range("G16")
$c.NumberFormat = "@"
$c.Value = "0.0333"
$c.Style = "Normal"
$ws.Range("H16").Value = 4
$ws.Range("A17").Value = 15
$c = $ws.Range("B17")
$c.NumberFormat = "@"
$c.Value = "011158"
$c.Style = "Normal"
$c = $ws.Range("C17")
$c.NumberFormat = "@"
$c.Value = "弘毅远方港股通智选领航混合C"
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.78"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "90.35"
$c.Style = "Normal"
$c = $ws.Range("F17")
$c.NumberFormat = "@"
$c.Value = "4.16"
$c.Style = "Normal"
$c = $ws.Range("G17")
$c.NumberFormat = "@"
$c.Value = "0.0324"
$c.Style = "Normal"
$ws.Range("H17").Value = 9
$ws.Range("A18").Value = 16
$c = $ws.Range("B18")
$c.NumberFormat = "@"
$c.Value = "011652"
$c.Style = "Normal"
$c = $ws.Range("C18")
$c.NumberFormat = "@"
$c.Value = "招商港股通核心精选股票C"
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.94"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "81.27"
$c.Style = "Normal"
$c = $ws.Range("F18")
$c.NumberFormat = "@"
$c.Value = "2.94"
$c.Style = "Normal"
$c = $ws.Range("G18")
$c.NumberFormat = "@"
$c.Value = "0.0276"
$c.Style = "Normal"
$ws.Range("H18").Value = 5
$ws.Range("A19").Value = 17
$c = $ws.Range("B19")
$c.NumberFormat = "@"
$c.Value = "011695"
$c.Style = "Normal"
$c = $ws.Range("C19")
$c.NumberFormat = "@"
$c.Value = "华泰紫金信息科技主题6个月定期开放混合C"
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.83"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "77.49"
$c.Style = "Normal"
$c = $ws.Range("F19")
$c.NumberFormat = "@"
$c.Value = "3.24"
$c.Style = "Normal"
$c = $ws.Range("G19")
$c.NumberFormat = "@"
$c.Value = "0.0269"
$c.Style = "Normal"
$ws.Range("H19").Value = 10
$ws.Range("A20").Value = 18
$c = $ws.Range("B20")
$c.NumberFormat = "@"
$c.Value = "013271"
$c.Style = "Normal"
$c = $ws.Range("C20")
$c.NumberFormat = "@"
$c.Value = "前海开源聚利一年持有混合C"
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.42"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "80.39"
$c.Style = "Normal"
$c = $ws.Range("F20")
$c.NumberFormat = "@"
$c.Value = "6.29"
$c.Style = "Normal"
$c = $ws.Range("G20")
$c.NumberFormat = "@"
$c.Value = "0.0264"
$c.Style = "Normal"
$ws.Range("H20").Value = 7
$ws.Range("A21").Value = 19
$c = $ws.Range("B21")
$c.NumberFormat = "@"
$c.Value = "006924"
$c.Style = "Normal"
$c = $ws.Range("C21")
$c.NumberFormat = "@"
$c.Value = "前海开源沪港深非周期性行业股票C"
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.22"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "93.77"
$c.Style = "Normal"
$c = $ws.Range("F21")
$c.NumberFormat = "@"
$c.Value = "6.16"
$c.Style = "Normal"
$c = $ws.Range("G21")
$c.NumberFormat = "@"
$c.Value = "0.0136"
$c.Style = "Normal"
$ws.Range("H21").Value = 4
$ws.Range("A22").Value = 20
$c = $ws.Range("B22")
$c.NumberFormat = "@"
$c.Value = "012315"
$c.Style = "Normal"
$c = $ws.Range("C22")
$c.NumberFormat = "@"
$c.Value = "创金合信港股通成长股票型发起式证券投资基金A"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.19"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "83.49"
$c.Style = "Normal"
$c = $ws.Range("F22")
$c.NumberFormat = "@"
$c.Value = "6.05"
$c.Style = "Normal"
$c = $ws.Range("G22")
$c.NumberFormat = "@"
$c.Value = "0.0115"
$c.Style = "Normal"
$ws.Range("H22").Value = 7
$ws.Range("A23").Value = 21
$c = $ws.Range("B23")
$c.NumberFormat = "@"
$c.Value = "010754"
$c.Style = "Normal"
$c = $ws.Range("C23")
$c.NumberFormat = "@"
$c.Value = "招商沪港深科技创新主题精选灵活配置混合C"
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.28"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "88.85"
$c.Style = "Normal"
$c = $ws.Range("F23")
$c.NumberFormat = "@"
$c.Value = "2.94"
$c.Style = "Normal"
$c = $ws.Range("G23")
$c.NumberFormat = "@"
$c.Value = "0.0082"
$c.Style = "Normal"
$ws.Range("H23").Value = 8
$ws.Range("A24").Value = 22
$c = $ws.Range("B24")
$c.NumberFormat = "@"
$c.Value = "012316"
$c.Style = "Normal"
$c = $ws.Range("C24")
$c.NumberFormat = "@"
$c.Value = "创金合信港股通成长股票型发起式证券投资基金C"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.10"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "83.49"
$c.Style = "Normal"
$c = $ws.Range("F24")
$c.NumberFormat = "@"
$c.Value = "6.05"
$c.Style = "Normal"
$c = $ws.Range("G24")
$c.NumberFormat = "@"
$c.Value = "0.0060"
$c.Style = "Normal"
$ws.Range("H24").Value = 7

# Clear the leftover template row (template had one more data row than needed)
$ws.Range("A25:H25").Clear()

# ================= Fill "总计" sheet =================
$ws = $newTotal
$ws.Range("B1").Value = "日期"
$ws.Range("C1").Value = "持有数量(只)"
$ws.Range("D1").Value = "持有市值(亿元)"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q1"
$ws.Range("C2").Value = 23
$ws.Range("D2").Value = 21.06
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "2021-Q4"
$ws.Range("C3").Value = 24
$ws.Range("D3").Value = 25.72
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "2021-Q3"
$ws.Range("C4").Value = 12
$ws.Range("D4").Value = 11.55
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "2021-Q2"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 0.17
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "2021-Q1"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 0.03
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "2020-Q4"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 0.17

# Clear leftover template columns/rows (template was an 8-col, 25-row sheet; totals sheet is 4-col, 7-row)
$ws.Range("E1:H25").Clear()
$ws.Range("A8:D25").Clear()

